$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 72
$ws.Range("A72").Value = 1739
$ws.Range("B72").Value = 102810
$ws.Range("C72").Value = "KLDT-O5WB"
$ws.Range("D72").Value = 54418
$ws.Range("E72").Value = "KLDT-E5WD"
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = "Bright flash x2 between 4b and LH antenna. Tungsten peak seem on Impurity signals on xpsedit"
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = "Two dots that do not move"
$ws.Range("K72").Value = "W"
$ws.Range("L72").Value = "Based on comment"
$ws.Range("M72").Value = "BEION4"
$ws.Range("N72").Value = "Single"

# New row 73
$ws.Range("A73").Value = 1568
$ws.Range("B73").Value = 103366
$ws.Range("C73").Value = "KLDT-O5WB"
$ws.Range("D73").Value = 51460
$ws.Range("E73").Value = "KLDT-E5WD"
$ws.Range("F73").Value = 54856
$ws.Range("G73").Value = "UFO's seen frame 335 from BEION plate and frame 336-337 from UDPT. Not seen on other cameras but increased radiation and Tungsten impurities found at same time on cview."
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = "Not seen in experimental cameras"
$ws.Range("K73").Value = "W"
$ws.Range("L73").Value = "Based on comment"
$ws.Range("M73").Value = "BEION4"
$ws.Range("N73").Value = "Single"

# New row 74
$ws.Range("A74").Value = 1730
$ws.Range("B74").Value = 102813
$ws.Range("C74").Value = "KLDT-O5WB"
$ws.Range("D74").Value = 4867101
$ws.Range("D74").NumberFormat = "#,##0"
$ws.Range("E74").Value = "KLDT-E5WD"
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = "3 bright spots in same location as previous pulse, still believed to be molybdenum"
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = "Taken from comment. It apears that with the Berilium is not able to track Mo down."
$ws.Range("J74").Value = "(213, 428)"
$ws.Range("K74").Value = "Mo"
$ws.Range("L74").Value = "Based on comment"
$ws.Range("L74").WrapText = $true
$ws.Range("M74").Value = "BEION4"
$ws.Range("N74").Value = "Single"

# New row 75
$ws.Range("A75").Value = 1733
$ws.Range("B75").Value = 102812
$ws.Range("C75").Value = "KLDT-O5WB"
$ws.Range("D75").Value = 5305262
$ws.Range("D75").NumberFormat = "#,##0"
$ws.Range("E75").Value = "KLDT-E5WD"
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = "potential molybdenum UFO as noted by spectroscopist"
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = "Literally same as above."
$ws.Range("J75").Value = "(236, 514)"
$ws.Range("K75").Value = "Mo"
$ws.Range("L75").Value = "Based on comment"
$ws.Range("L75").WrapText = $true
$ws.Range("M75").Value = "BEION4"
$ws.Range("N75").Value = "Single"
